# Apply updated crypto price/volume figures to the worksheet (row 2 = Bitcoin ... row 51 = EnergySwap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.633.10'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '2.621.40'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'" + '594.94'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").Value = "'" + '150.09'
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").Value = "'" + '5.71'
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("E11").Value = '  +2.88%  '
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").Value = "'" + '27.66'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").Value = '3.093.76'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '63.451.70'
$ws.Range("E15").Value = '  -0.01%  '
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").Value = '2.641.11'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").Value = "'" + '12.33'
$ws.Range("E18").Value = '  +7.19%  '
$ws.Range("D19").Value = "'" + '4.66'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").Value = "'" + '347.30'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = "'" + '0.997'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").Value = "'" + '5.72'
$ws.Range("E23").Value = '  +2.36%  '
$ws.Range("D24").Value = "'" + '66.42'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("E25").Value = '  +10.53%  '
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = "'" + '9.22'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").Value = "'" + '576.18'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").Value = "'" + '8.23'
$ws.Range("E29").Value = '  +3.11%  '
$ws.Range("D30").Value = "'" + '0.163'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = "'" + '2.05'
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("E33").Value = '  +2.72%  '
$ws.Range("E34").Value = '  -0.27%  '
$ws.Range("D35").Value = "'" + '5.27'
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").Value = "'" + '168.69'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").Value = "'" + '0.999'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").Value = "'" + '19.36'
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("D42").Value = "'" + '168.39'
$ws.Range("E42").Value = '  -0.49%  '
$ws.Range("D43").Value = "'" + '39.88'
$ws.Range("E44").Value = '  +4.08%  '
$ws.Range("E45").Value = '  +5.14%  '
$ws.Range("D46").Value = "'" + '21.42'
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("D47").Value = "'" + '0.627'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("E49").Value = '  +4.80%  '
$ws.Range("D50").Value = "'" + '0.0962'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").Value = "'" + '19.30'
$ws.Range("E51").Value = '  +2.43%  '
